$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.823277354240417
$ws.Range("B1").Value = 2.302476406097412
$ws.Range("C1").Value = 2.434069395065308
$ws.Range("D1").Value = 2.927975177764893
$ws.Range("E1").Value = 2.070175170898438
